$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old totals row (row 6), pushing the totals
# row down to 7 and the "Final Grade" row down to 11.
$ws.Rows("6:6").Insert(-4121)

# The bottom border/fill combo that used to sit under "hw4?" (row 5) now
# belongs under the new last homework row (row 6) - copy that formatting
# down to the new row before changing row 5's own formatting.
$ws.Range("B5:C5").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "hw5" homework row (row 6).
$ws.Range("A6").Value = "hw5"
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 100

# Row 5 ("hw4?") no longer has the divider border; C5 keeps only its fill.
$ws.Range("B5").Borders.Item(9).LineStyle = -4142
$ws.Range("C5").Borders.Item(9).LineStyle = -4142

# Update the totals/summary formulas (now on row 7) to include the new
# homework row.
$ws.Range("B7").Formula = "=SUM(B2:B6)"
$ws.Range("C7").Formula = "=SUM(C2:C6)"
$ws.Range("E7").Formula = "=(C7/B7)*D7"
$ws.Range("G7").Formula = "=SUM(G2:G6)"
$ws.Range("H7").Formula = "=SUM(H2:H6)"
$ws.Range("J7").Formula = "=(H7/G7)*I7"

# Final grade row (now row 11).
$ws.Range("B11").Formula = "=E7+J7"

# Move the active selection to A7, matching the post-edit cursor position.
$ws.Range("A7").Select()
